# Excel COM-interop script: Add new Marketing ROI and Marketing Efficiency
# metrics to the YOY table (inserted as new rows 9 & 10), shifting the
# existing rows down by two. Also renames the "AFC / RCM Payments"
# responsibility label to "AFC / Vendor Payments", and corrects the
# August actual expense figure (F8) which feeds the new ratios.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert two blank rows above the old row 9 (old rows 9-44 -> 11-46)
# ---------------------------------------------------------------------
$ws.Range("A9:A10").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2) Correct the August 2025 actual expense value that the new ratios
#    depend on (was a placeholder estimate, now the actual figure).
# ---------------------------------------------------------------------
$ws.Range("F8").Value = 66195

# ---------------------------------------------------------------------
# 3) New row 9: MARKETING_ROI / "Marketing ROI (1$ Earned…)"
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "MARKETING_ROI"
$ws.Range("B9").Value = "Marketing ROI (1$ Earned…)"
$ws.Range("C9").Value = "YOY Expense & Profitability Analysis"
$ws.Range("D9").Value = "Owner Controlled"
$ws.Range("E9").Formula = "=((E3-E8)/E8)"
$ws.Range("F9").Formula = "=(F3-F8)/F8"
$ws.Range("G9").Formula = "=(F9-E9)/E9"
$ws.Range("H9").Value = 26
$ws.Range("E9:F9").Style = "Currency"

# ---------------------------------------------------------------------
# 4) New row 10: MARKETING_EFFICIENCY / "Marketing Efficiency (# of
#    visits per dollar)"
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "MARKETING_EFFICIENCY"
$ws.Range("B10").Value = "Marketing Efficiency (# of visits per dollar)"
$ws.Range("C10").Value = "YOY Expense & Profitability Analysis"
$ws.Range("D10").Value = "Owner Controlled"
$ws.Range("E10").Formula = "=E4/E8"
$ws.Range("F10").Formula = "=F4/F8"
$ws.Range("G10").Formula = "=(F10-E10)/E10"
$ws.Range("H10").Value = 9
$ws.Range("E10:F10").Style = "Comma"

# ---------------------------------------------------------------------
# 5) Row 11 (old row 9, "Revenue Proportional"): the Responsibility
#    label moves from "AFC / RCM Payments" to "AFC / Vendor Payments"
# ---------------------------------------------------------------------
$ws.Range("D11").Value = "AFC / Vendor Payments"

# ---------------------------------------------------------------------
# 6) Update the selection / view to match the edited area
# ---------------------------------------------------------------------
$ws.Range("D26").Select()

$wb.Save()
